# Generate Report for Handback
# Adds a new handback row (aff6cbec-58b3-460e-b055-9fc8d82fe812.md) to the
# Overview, zh-cn and de-de sheets/tables.

$wb = $excel.ActiveWorkbook

$guidFile   = "aff6cbec-58b3-460e-b055-9fc8d82fe812.md"
$guidPath   = "e2e\aff6cbec-58b3-460e-b055-9fc8d82fe812.md"
$statusSync = "Handed back: in sync with en-US"

$zhXlf = "aff6cbec-58b3-460e-b055-9fc8d82fe812.095c20b14c39c7a7cdcffbe903e0ced8b1560904.zh-cn.xlf"
$deXlf = "aff6cbec-58b3-460e-b055-9fc8d82fe812.095c20b14c39c7a7cdcffbe903e0ced8b1560904.de-de.xlf"

$zhHoDate = "2016-08-19 02:41:58"
$zhHbDate = "2016-08-19 02:42:28"
$deHoDate = "2016-08-19 02:42:09"
$deHbDate = "2016-08-19 02:42:35"

# ---------------------------------------------------------------------------
# Overview sheet (row 4)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A4").Value = $guidFile
$wsOverview.Range("B4").Value = $guidPath
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("E4").Value = $statusSync
$wsOverview.Range("F4").Value = $statusSync
$wsOverview.Range("G4").Value = $deHoDate

$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aff6cbec58b3460eb0559fc8d82fe812e812aff6/$guidPath", "", "", $guidPath) | Out-Null
$wsOverview.Range("B4").Font.Underline = $true
$wsOverview.Range("B4").Font.Color = 15570276

# ---------------------------------------------------------------------------
# zh-cn sheet (row 4)
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null

$wsZh.Range("A4").Value = $guidFile
$wsZh.Range("B4").Value = ".md"
$wsZh.Range("C4").Value = $statusSync
$wsZh.Range("D4").Value = "e2e"
$wsZh.Range("E4").Value = "ht"
$wsZh.Range("F4").Value = "True"
$wsZh.Range("G4").Value = $zhXlf
$wsZh.Range("H4").Value = $zhHoDate
$wsZh.Range("I4").Value = $guidFile
$wsZh.Range("J4").Value = $zhXlf
$wsZh.Range("K4").Value = $zhHbDate
$wsZh.Range("L4").Value = ""
$wsZh.Range("M4").Value = "True"
$wsZh.Range("N4").Value = ""
$wsZh.Range("O4").Value = "False"
$wsZh.Range("P4").Value = ""

$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aff6cbec58b3460eb0559fc8d82fe812e812aff6/e2e/$guidFile", "", "", $guidFile) | Out-Null
$wsZh.Range("A4").Font.Underline = $true
$wsZh.Range("A4").Font.Color = 15570276

$wsZh.Hyperlinks.Add($wsZh.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/aff6cbec58b3460eb0559fc8d82fe812e812aff6/e2e/$guidFile", "", "", $guidFile) | Out-Null
$wsZh.Range("I4").Font.Underline = $true
$wsZh.Range("I4").Font.Color = 15570276

$wsZh.Range("H4").NumberFormat = $wsZh.Range("H2").NumberFormat
$wsZh.Range("K4").NumberFormat = $wsZh.Range("K2").NumberFormat

# ---------------------------------------------------------------------------
# de-de sheet (row 4)
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null

$wsDe.Range("A4").Value = $guidFile
$wsDe.Range("B4").Value = ".md"
$wsDe.Range("C4").Value = $statusSync
$wsDe.Range("D4").Value = "e2e"
$wsDe.Range("E4").Value = "ht"
$wsDe.Range("F4").Value = "True"
$wsDe.Range("G4").Value = $deXlf
$wsDe.Range("H4").Value = $deHoDate
$wsDe.Range("I4").Value = $guidFile
$wsDe.Range("J4").Value = $deXlf
$wsDe.Range("K4").Value = $deHbDate
$wsDe.Range("L4").Value = ""
$wsDe.Range("M4").Value = "True"
$wsDe.Range("N4").Value = ""
$wsDe.Range("O4").Value = "False"
$wsDe.Range("P4").Value = ""

$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aff6cbec58b3460eb0559fc8d82fe812e812aff6/e2e/$guidFile", "", "", $guidFile) | Out-Null
$wsDe.Range("A4").Font.Underline = $true
$wsDe.Range("A4").Font.Color = 15570276

$wsDe.Hyperlinks.Add($wsDe.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/aff6cbec58b3460eb0559fc8d82fe812e812aff6/e2e/$guidFile", "", "", $guidFile) | Out-Null
$wsDe.Range("I4").Font.Underline = $true
$wsDe.Range("I4").Font.Color = 15570276

$wsDe.Range("H4").NumberFormat = $wsDe.Range("H2").NumberFormat
$wsDe.Range("K4").NumberFormat = $wsDe.Range("K2").NumberFormat
